$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.346.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.841.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.610"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.18%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.70"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +11.67%  "
$ws.Range("E9").Value = "  +6.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0693"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("E11").Value = "  +3.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.109.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.844.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.669"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.383.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0799"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "243.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.97%  "
$ws.Range("E22").Value = "  +13.15%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("E29").Value = "  +13.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.412.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +40.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0545"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.11%  "
$ws.Range("E33").Value = "  +4.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "95.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +15.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.683"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.87%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.340.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("E40").Value = "  +5.92%  "
$ws.Range("E41").Value = "  +3.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.03%  "
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("E47").Value = "  +8.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0519"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.010.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.06%  "
